$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 750.7308
$ws.Range("I125").Value = 503
$ws.Range("J125").Value = 963.0714
$ws.Range("K125").Value = 4527
$ws.Range("L125").Value = 8667.642600000001
$ws.Range("M125").Value = -2067
$ws.Range("N125").Value = -13587.6426
$ws.Range("H137").Value = 1881.3636
$ws.Range("I137").Value = 1390.3636
$ws.Range("J137").Value = 2372.3635
$ws.Range("K137").Value = 4171.0908
$ws.Range("L137").Value = 7117.0905
$ws.Range("M137").Value = -1621.0908
$ws.Range("N137").Value = -12217.0905
$ws.Range("H138").Value = 1473.7709
$ws.Range("J138").Value = 1855.3541
$ws.Range("L138").Value = 5566.0623
$ws.Range("N138").Value = -15846.0623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1908.1919
$ws.Range("I32").Value = 1951.2604
$ws.Range("J32").Value = 530
$ws.Range("K32").Value = 1951.2604
$ws.Range("L32").Value = 530
$ws.Range("M32").Value = -1664.2604
$ws.Range("N32").Value = -1104
$ws.Range("H61").Value = 1305.5652
$ws.Range("I61").Value = 917.5714
$ws.Range("J61").Value = 2256.15
$ws.Range("K61").Value = 917.5714
$ws.Range("L61").Value = 2256.15
$ws.Range("M61").Value = -705.5714
$ws.Range("N61").Value = -2680.15
$ws.Range("H74").Value = 902.8293
$ws.Range("I74").Value = 874
$ws.Range("J74").Value = 1042.8572
$ws.Range("K74").Value = 874
$ws.Range("L74").Value = 1042.8572
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = -2790.8572
$ws.Range("H77").Value = 902.8293
$ws.Range("I77").Value = 874
$ws.Range("J77").Value = 1042.8572
$ws.Range("K77").Value = 4370
$ws.Range("L77").Value = 5214.286
$ws.Range("M77").Value = -2
$ws.Range("N77").Value = -13950.286
$ws.Range("H132").Value = 3532.7896
$ws.Range("I132").Value = 3644.818
$ws.Range("J132").Value = 3153.6155
$ws.Range("K132").Value = 10934.454
$ws.Range("L132").Value = 9460.8465
$ws.Range("M132").Value = -8404.454000000002
$ws.Range("N132").Value = -14520.8465
$ws.Range("H136").Value = 1305.5652
$ws.Range("I136").Value = 917.5714
$ws.Range("J136").Value = 2256.15
$ws.Range("K136").Value = 2752.7142
$ws.Range("L136").Value = 6768.450000000001
$ws.Range("M136").Value = -202.7142000000003
$ws.Range("N136").Value = -11868.45

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1927.6097
$ws.Range("I134").Value = 1644.7333
$ws.Range("J134").Value = 2699.0908
$ws.Range("K134").Value = 4934.199900000001
$ws.Range("L134").Value = 8097.2724
$ws.Range("M134").Value = -2399.199900000001
$ws.Range("N134").Value = -13167.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 737.5263
$ws.Range("I16").Value = 680.2353000000001
$ws.Range("J16").Value = 1224.5
$ws.Range("K16").Value = 680.2353000000001
$ws.Range("L16").Value = 1224.5
$ws.Range("M16").Value = -393.2353000000001
$ws.Range("N16").Value = -1798.5
$ws.Range("H31").Value = 2875688
$ws.Range("I31").Value = 1541.125
$ws.Range("J31").Value = 9262681
$ws.Range("K31").Value = 1541.125
$ws.Range("L31").Value = 9262681
$ws.Range("M31").Value = -1246.125
$ws.Range("N31").Value = -9263271
$ws.Range("H34").Value = 2875688
$ws.Range("I34").Value = 1541.125
$ws.Range("J34").Value = 9262681
$ws.Range("K34").Value = 1541.125
$ws.Range("L34").Value = 9262681
$ws.Range("M34").Value = -1339.125
$ws.Range("N34").Value = -9263085
$ws.Range("H86").Value = 90911030
$ws.Range("I86").Value = 125001760
$ws.Range("J86").Value = 2433.3333
$ws.Range("K86").Value = 125001760
$ws.Range("L86").Value = 2433.3333
$ws.Range("M86").Value = -125000637
$ws.Range("N86").Value = -4679.3333
$ws.Range("H89").Value = 90911030
$ws.Range("I89").Value = 125001760
$ws.Range("J89").Value = 2433.3333
$ws.Range("K89").Value = 625008800
$ws.Range("L89").Value = 12166.6665
$ws.Range("M89").Value = -625003184
$ws.Range("N89").Value = -23398.6665
$ws.Range("H113").Value = 737.5263
$ws.Range("I113").Value = 680.2353000000001
$ws.Range("J113").Value = 1224.5
$ws.Range("K113").Value = 680.2353000000001
$ws.Range("L113").Value = 1224.5
$ws.Range("M113").Value = 1489.7647
$ws.Range("N113").Value = -5564.5
$ws.Range("H134").Value = 1434.2632
$ws.Range("I134").Value = 1326.9697
$ws.Range("J134").Value = 2142.4
$ws.Range("K134").Value = 3980.9091
$ws.Range("L134").Value = 6427.200000000001
$ws.Range("M134").Value = -1445.9091
$ws.Range("N134").Value = -11497.2
$ws.Range("H135").Value = 26142.857
$ws.Range("J135").Value = 26142.857
$ws.Range("L135").Value = 26142.857
$ws.Range("N135").Value = -36282.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1206.5
$ws.Range("I5").Value = 274.44446
$ws.Range("J5").Value = 2255.0625
$ws.Range("K5").Value = 823.33338
$ws.Range("L5").Value = 6765.1875
$ws.Range("M5").Value = -711.33338
$ws.Range("N5").Value = -6989.1875
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("H97").Value = 276.55
$ws.Range("I97").Value = 280
$ws.Range("J97").Value = 275.07144
$ws.Range("K97").Value = 840
$ws.Range("L97").Value = 825.21432
$ws.Range("M97").Value = -344
$ws.Range("N97").Value = -1817.21432
$ws.Range("H122").Value = 1404.1666
$ws.Range("I122").Value = 1281.091
$ws.Range("J122").Value = 1597.5714
$ws.Range("K122").Value = 11529.819
$ws.Range("L122").Value = 14378.1426
$ws.Range("M122").Value = -9079.819
$ws.Range("N122").Value = -19278.1426
$ws.Range("H131").Value = 5632049
$ws.Range("I131").Value = 11530
$ws.Range("J131").Value = 8793591
$ws.Range("K131").Value = 34590
$ws.Range("L131").Value = 26380773
$ws.Range("M131").Value = -29550
$ws.Range("N131").Value = -26390853
$ws.Range("H133").Value = 2766.25
$ws.Range("I133").Value = 2855
$ws.Range("K133").Value = 8565
$ws.Range("M133").Value = -3505
$ws.Range("H134").Value = 1891.7059
$ws.Range("I134").Value = 1638.25
$ws.Range("K134").Value = 4914.75
$ws.Range("M134").Value = 155.25
$ws.Range("H135").Value = 1206.5
$ws.Range("I135").Value = 274.44446
$ws.Range("J135").Value = 2255.0625
$ws.Range("K135").Value = 2470.00014
$ws.Range("L135").Value = 20295.5625
$ws.Range("M135").Value = 64.9998599999999
$ws.Range("N135").Value = -25365.5625
$ws.Range("H136").Value = 1283.3
$ws.Range("J136").Value = 3155.3333
$ws.Range("L136").Value = 9465.999899999999
$ws.Range("N136").Value = -19665.9999
$ws.Range("H137").Value = 66206860
$ws.Range("I137").Value = 30314084
$ws.Range("J137").Value = 105688910
$ws.Range("K137").Value = 90942252
$ws.Range("L137").Value = 317066730
$ws.Range("M137").Value = -90937152
$ws.Range("N137").Value = -317076930
$ws.Range("H138").Value = 1714.7059
$ws.Range("I138").Value = 1516.6666
$ws.Range("J138").Value = 3200
$ws.Range("K138").Value = 4549.9998
$ws.Range("L138").Value = 9600
$ws.Range("M138").Value = 590.0002000000004
$ws.Range("N138").Value = -19880
$ws.Range("H139").Value = 2756.3635
$ws.Range("I139").Value = 1990
$ws.Range("K139").Value = 5970
$ws.Range("M139").Value = -830
$ws.Range("H140").Value = 1186.0344
$ws.Range("I140").Value = 1067.6786
$ws.Range("J140").Value = 4500
$ws.Range("K140").Value = 3203.0358
$ws.Range("L140").Value = 13500
$ws.Range("M140").Value = 1976.9642
$ws.Range("N140").Value = -23860
$ws.Range("H141").Value = 2336.5356
$ws.Range("I141").Value = 2225.9583
$ws.Range("K141").Value = 6677.874899999999
$ws.Range("M141").Value = -1497.874899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 26601.666
$ws.Range("I14").Value = 25000
$ws.Range("J14").Value = 29805
$ws.Range("K14").Value = 25000
$ws.Range("L14").Value = 29805
$ws.Range("M14").Value = -24832
$ws.Range("N14").Value = -30141
$ws.Range("H70").Value = 9811675
$ws.Range("I70").Value = 11864451
$ws.Range("J70").Value = 3966.6667
$ws.Range("K70").Value = 11864451
$ws.Range("L70").Value = 3966.6667
$ws.Range("M70").Value = -11864181
$ws.Range("N70").Value = -4506.6667
$ws.Range("H73").Value = 9811675
$ws.Range("I73").Value = 11864451
$ws.Range("J73").Value = 3966.6667
$ws.Range("K73").Value = 11864451
$ws.Range("L73").Value = 3966.6667
$ws.Range("M73").Value = -11863515
$ws.Range("N73").Value = -5838.6667
$ws.Range("H107").Value = 325.5
$ws.Range("I107").Value = 400.5
$ws.Range("J107").Value = 263
$ws.Range("K107").Value = 400.5
$ws.Range("L107").Value = 263
$ws.Range("M107").Value = 1519.5
$ws.Range("N107").Value = -4103

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2242.7532
$ws.Range("I136").Value = 2434.2363
$ws.Range("J136").Value = 1764.0454
$ws.Range("K136").Value = 7302.7089
$ws.Range("L136").Value = 5292.1362
$ws.Range("M136").Value = -4752.7089
$ws.Range("N136").Value = -10392.1362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1395.2059
$ws.Range("I132").Value = 714.381
$ws.Range("J132").Value = 2495
$ws.Range("K132").Value = 2143.143
$ws.Range("L132").Value = 7485
$ws.Range("M132").Value = 386.857
$ws.Range("N132").Value = -12545
$ws.Range("H136").Value = 2218.8735
$ws.Range("I136").Value = 2451.9524
$ws.Range("J136").Value = 1301.125
$ws.Range("K136").Value = 7355.8572
$ws.Range("L136").Value = 3903.375
$ws.Range("M136").Value = -4805.8572
$ws.Range("N136").Value = -9003.375
